$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: mark E10 as "Paid"
$ws.Range("E10").Value = "Paid"

# Row 11: mark E11 as "Paid"
$ws.Range("E11").Value = "Paid"

# Row 12: fill in payment date and mark as Paid
$ws.Range("A12").Value2 = 45234
$ws.Range("A12").NumberFormat = "m/d/yy"
$ws.Range("E12").Value = "Paid"

# Row 13: fill in payment date, amount formula, and mode of payment
$ws.Range("A13").Value2 = 45236
$ws.Range("A13").NumberFormat = "m/d/yy"
$ws.Range("B13").Formula = "=2700+900"
$ws.Range("C13").Value = "NEFT"

# Update the active cell selection
$ws.Range("H25").Select()
